$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D6").Value = -7.621499999999995
$ws.Range("C7").Value = -13.0204
$ws.Range("A9").Value = -21.73120000000001
$ws.Range("C12").Value = -10.8106
$ws.Range("C14").Value = -13.42709999999999
$ws.Range("D15").Value = -8.900299999999998
$ws.Range("A18").Value = -22.17930000000002
$ws.Range("A20").Value = -21.23709999999999
$ws.Range("C26").Value = -12.5696
$ws.Range("A27").Value = -21.8994
$ws.Range("C27").Value = -12.9353
$ws.Range("C29").Value = -11.1005
$ws.Range("D33").Value = -7.743999999999999
$ws.Range("A35").Value = -19.95789999999998
$ws.Range("D35").Value = -7.943199999999995
$ws.Range("C37").Value = -13.98689999999998
$ws.Range("C38").Value = -13.5157
$ws.Range("D38").Value = -8.681499999999991
$ws.Range("D43").Value = -8.141500000000004
$ws.Range("D44").Value = -7.801099999999998
$ws.Range("D47").Value = -7.484800000000003
$ws.Range("C51").Value = -12.3072
$ws.Range("D51").Value = -7.786700000000002
$ws.Range("C52").Value = -11.2789
$ws.Range("C55").Value = -13.5659
$ws.Range("D57").Value = -8.261500000000002
$ws.Range("D63").Value = -8.020900000000001
$ws.Range("A69").Value = -21.72310000000001
$ws.Range("C69").Value = -11.1274
$ws.Range("C70").Value = -12.76300000000001
$ws.Range("D70").Value = -8.216199999999995
$ws.Range("A76").Value = -19.95679999999999
$ws.Range("A78").Value = -20.01129999999998
$ws.Range("C81").Value = -13.4148
$ws.Range("A82").Value = -22.03930000000001
$ws.Range("A83").Value = -21.6314
$ws.Range("C83").Value = -12.0255
$ws.Range("D88").Value = -8.077099999999998
$ws.Range("A93").Value = -21.0356
$ws.Range("D99").Value = -7.752099999999993
$ws.Range("C102").Value = -13.3936
